{"js": "// The document contains three quiz questions, each with a title paragraph\n// followed by four answer-choice paragraphs (the italic choice is the\n// correct answer). This edit rewrites question 1 and question 3 in place:\n// question 1's title/choices become the old question 3's content, and\n// question 3's title/choices become the old question 1's content, with the\n// answer choices also reshuffled (the italic/correct marker moves with its\n// matching text). Question 2 keeps its title but its four choices are\n// likewise reshuffled (still ending with RIP as the correct/italic choice).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nparagraphs.items.forEach((p) => p.font.load(\"italic\"));\nawait context.sync();\n\n// Target state for paragraph index -> { text, italic } covering the whole\n// body (index 0 is the empty heading paragraph, which is left untouched).\nconst targets = [\n  null, // 0: heading paragraph - unchanged\n  { text: \"1 Which of the following devices is responsible for selecting the best path for a datagram?\", italic: false },\n  { text: \"NIC\", italic: false },\n  { text: \"Switch\", italic: false },\n  { text: \"Router\", italic: true },\n  { text: \"Hub\", italic: false },\n  { text: \"2 Which of the following protocols use distance-vector routing?\", italic: false },\n  { text: \"DHCP\", italic: false },\n  { text: \"OSPF\", italic: false },\n  { text: \"BGP\", italic: false },\n  { text: \"RIP\", italic: true },\n  { text: \"3 A socket is composed of?\", italic: false },\n  { text: \"MAC Address\", italic: false },\n  { text: \"Port Number\", italic: true },\n  { text: \"IP Address\", italic: true },\n  { text: \"URL\", italic: false },\n];\n\nconst count = Math.min(paragraphs.items.length, targets.length);\nfor (let i = 0; i < count; i++) {\n  const target = targets[i];\n  if (!target) continue;\n\n  const paragraph = paragraphs.items[i];\n  const textChanged = paragraph.text !== target.text;\n  const italicChanged = paragraph.font.italic !== target.italic;\n  if (!textChanged && !italicChanged) continue;\n\n  if (textChanged) {\n    // Replace the paragraph's text while keeping a single run.\n    const range = paragraph.getRange();\n    range.insertText(target.text, Word.InsertLocation.replace);\n  }\n\n  if (italicChanged) {\n    paragraph.font.italic = target.italic;\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains three quiz questions, each with a title paragraph\n# followed by four answer-choice paragraphs (the italic choice is the\n# correct answer). This edit rewrites question 1 and question 3 in place:\n# question 1's title/choices become the old question 3's content, and\n# question 3's title/choices become the old question 1's content, with the\n# answer choices also reshuffled (the italic/correct marker moves with its\n# matching text). Question 2 keeps its title but its four choices are\n# likewise reshuffled (still ending with RIP as the correct/italic choice).\n\n$d = $word.ActiveDocument\n\n# Target state keyed by 1-based paragraph index (Paragraphs(1) is the empty\n# heading paragraph, left untouched).\n$targets = @{\n    2  = @{ Text = \"1 Which of the following devices is responsible for selecting the best path for a datagram?\"; Italic = $false }\n    3  = @{ Text = \"NIC\"; Italic = $false }\n    4  = @{ Text = \"Switch\"; Italic = $false }\n    5  = @{ Text = \"Router\"; Italic = $true }\n    6  = @{ Text = \"Hub\"; Italic = $false }\n    7  = @{ Text = \"2 Which of the following protocols use distance-vector routing?\"; Italic = $false }\n    8  = @{ Text = \"DHCP\"; Italic = $false }\n    9  = @{ Text = \"OSPF\"; Italic = $false }\n    10 = @{ Text = \"BGP\"; Italic = $false }\n    11 = @{ Text = \"RIP\"; Italic = $true }\n    12 = @{ Text = \"3 A socket is composed of?\"; Italic = $false }\n    13 = @{ Text = \"MAC Address\"; Italic = $false }\n    14 = @{ Text = \"Port Number\"; Italic = $true }\n    15 = @{ Text = \"IP Address\"; Italic = $true }\n    16 = @{ Text = \"URL\"; Italic = $false }\n}\n\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if (-not $targets.ContainsKey($i)) { continue }\n\n    $target = $targets[$i]\n    $p = $d.Paragraphs($i)\n    $r = $p.Range\n\n    if ($r.Text -ne $target.Text) {\n        $r.Text = $target.Text\n    }\n\n    $wantItalic = 0\n    if ($target.Italic) { $wantItalic = -1 }\n    if ($r.Italic -ne $wantItalic) {\n        $r.Italic = $wantItalic\n    }\n}\n"}
